# Ikura - Project Tracker.docx edit
# 1) "Make fourth page (login and signup as" / " one??)" were two separate runs split by
#    the _GoBack bookmark. Merge them into a single run with the full text, which also
#    drops the now-stale bookmark from this spot (it gets re-added further down, after the
#    new "To do:" paragraph, as part of the appended OOXML below).
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Make fourth page (login and signup as one??)", $false, $false, $false, $false, $false,
    $true, 1, $false, "Make fourth page (login and signup as one??)", 2) | Out-Null

# 2) Append the new "Week 1 - Day 1" status-update section (Done / Working on / To do, with
#    nested sub-bullets) plus the trailing blank paragraph, as raw WordprocessingML inserted
#    right after the last paragraph in the document.
$xmlFragment = @'
<w:p/><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Week 1</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> – Day 1</w:t></w:r></w:p><w:p><w:r><w:t>Done:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Make fourth page (login and signup as one??)</w:t></w:r><w:r><w:t xml:space="preserve"> – made whiteboard version</w:t></w:r></w:p><w:p><w:r><w:t>Working on:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Creating a skeleton of:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Base html (with minimal bootstrap on base)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>All html pages w/ connections to base html</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Routes from one html page to another</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Creating server.py doc</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Working on drawing out relationships between everything</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Make sure that my database structure is scalable. I want to be able to use the three (money, time, sanity) factors and also maybe integrate mint’s </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>api</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> for later projects. Also have to consider how </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>twilio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> app text updates with monthly payments will affect my database. </w:t></w:r></w:p><w:p><w:r><w:t>To do:</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Math excel spread sheet</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Start with one card</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Move up to two</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Then three</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Start thinking about how best to store all pieces of data for each step of the calculations</w:t></w:r></w:p><w:p/>
'@

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $xmlFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$endPos = $d.Content.End
$insertionRange = $d.Range($endPos, $endPos)
$insertionRange.InsertXML($packageXml)

Write-Output "Applied Week 1 - Day 1 update"
